$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 44 with the new activity log entry
$ws.Range("B44").Value = 6977
$ws.Range("C44").Value = 43926
$ws.Range("D44").Value = 0.068749999999999992
$ws.Range("E44").Value = 0.074305555555555555
$ws.Range("G44").Value = "Added waveforms for LogicUnit.vhd to document"

# Update the active selection to A44, matching the saved view state
$ws.Range("A44").Select()
